$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B41").Value = "ESCs"
